$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# Rows 3-41 in column D ("Runmode") switch from "Y" to "N" (row 2 stays "Y").
$ws.Range("D3:D41").Value = "N"

# Scroll the window so row 10 is at the top, and select D3:D41 with the
# active cell on D3 (matches the saved selection/topLeftCell in the sheet view).
$ws.Range("D3:D41").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
